$d = $word.ActiveDocument

$replacements = @(
    @("2026-02-12 Thursday", "2026-02-13 Friday"),
    @("27×79=2133", "91×35=3185"),
    @("56×38=2128", "22×33=726"),
    @("28×84=2352", "39×20=780"),
    @("65×59=3835", "30×65=1950"),
    @("87×32=2784", "62×23=1426"),
    @("11×98=1078", "38×95=3610"),
    @("99×26=2574", "78×91=7098"),
    @("24×49=1176", "13×36=468"),
    @("28×68=1904", "93×16=1488"),
    @("14×91=1274", "62×31=1922"),
    @("29×64=1856", "29×88=2552"),
    @("67×53=3551", "29×27=783"),
    @("30×16=480", "85×54=4590"),
    @("55×11=605", "16×73=1168"),
    @("68×83=5644", "88×41=3608"),
    @("24×11=264", "81×73=5913"),
    @("98×64=6272", "13×39=507"),
    @("11×89=979", "90×14=1260"),
    @("52×16=832", "80×61=4880"),
    @("56×23=1288", "21×36=756"),
    @("72×74=5328", "44×66=2904"),
    @("70×34=2380", "74×41=3034"),
    @("19×57=1083", "12×25=300"),
    @("28×37=1036", "65×75=4875"),
    @("77×47=3619", "14×14=196")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $new, 2)
}
